# Fix typo in police_avoidance3 question text: drop the stray
# "Plan drive/walk " prefix that was accidentally left in the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "How often do you plan where to drive/ walk to make sure that you do not attract attention from the police?"

# Move the active selection to B5 (matches the saved cursor position).
$ws.Range("B5").Select()
